# fitur : hapus data parkir selesai
# Remove the completed parking entry (row 3 - "SN12XKBEG18L"), shifting all
# subsequent rows up by one. Then update the (new) last row with the
# freshly-arrived, still-parked vehicle's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the completed parking record in row 3; rows below shift up.
$ws.Rows(3).Delete()

# After the shift, row 7 is now the last data row and represents a brand
# new, still in-progress entry (no exit time / duration / cost yet).
$ws.Cells.Item(7, 1).Value = "TXDSWY5EZFB6"
$ws.Cells.Item(7, 2).Value = "D 0423 JS"
$ws.Cells.Item(7, 3).Value = " "
$ws.Cells.Item(7, 4).Value = "2025-02-03 18:56:00"
$ws.Cells.Item(7, 5).Value = " "
$ws.Cells.Item(7, 6).Value = " "
$ws.Cells.Item(7, 7).Value = " "
$ws.Cells.Item(7, 8).Value = " "
$ws.Cells.Item(7, 9).Value = "Reza Ramdan Permana"
$ws.Cells.Item(7, 10).Value = "./capture/masuk/TXDSWY5EZFB6.png"
$ws.Cells.Item(7, 11).Value = " "
